$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all existing data rows
# (2 through 338) from 45172 (2023-09-03) to 45175 (2023-09-06).
$ws.Range("C2:C338").Value = 45175

# Row 338 gains an explicit row height (ht="15" customHeight="1") in the
# target file; setting RowHeight to its current value reproduces that.
$ws.Rows.Item(338).RowHeight = 15

# Append the new record as row 339.
$ws.Cells.Item(339, 1).Value = "A 41205-2023"

$ws.Cells.Item(339, 2).Value = 45174
$ws.Cells.Item(339, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(339, 3).Value = 45175
$ws.Cells.Item(339, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(339, 4).Value = "ÖREBRO LÄN"
$ws.Cells.Item(339, 5).Value = "HÄLLEFORS"

$ws.Cells.Item(339, 7).Value = 4
$ws.Cells.Item(339, 8).Value = 0
$ws.Cells.Item(339, 9).Value = 0
$ws.Cells.Item(339, 10).Value = 0
$ws.Cells.Item(339, 11).Value = 0
$ws.Cells.Item(339, 12).Value = 0
$ws.Cells.Item(339, 13).Value = 0
$ws.Cells.Item(339, 14).Value = 0
$ws.Cells.Item(339, 15).Value = 0
$ws.Cells.Item(339, 16).Value = 0
$ws.Cells.Item(339, 17).Value = 0

# Column R carries a wrap-text style even though it stays empty.
$ws.Cells.Item(339, 18).WrapText = $true
